$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for Price (D) column cells whose new value would
# otherwise be auto-parsed by Excel as a number (losing literal formatting
# like trailing zeros, e.g. "1.00" -> 1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.873.68"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.733.00"
$ws.Range("E3").Value = "  +3.32%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "602.73"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D6").Value = "164.51"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "2.731.49"
$ws.Range("E9").Value = "  +3.26%  "

$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("D13").Value = "5.31"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").Value = "28.65"
$ws.Range("E14").Value = "  +2.11%  "

$ws.Range("D15").Value = "3.229.16"
$ws.Range("E15").Value = "  +3.30%  "

$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("D17").Value = "68.875.53"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "2.716.54"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").Value = "11.90"
$ws.Range("E19").Value = "  +4.45%  "

$ws.Range("E20").Value = "  +5.62%  "

$ws.Range("D21").Value = "366.92"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").Value = "4.55"
$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  +3.19%  "

$ws.Range("D25").Value = "73.91"
$ws.Range("E25").Value = "  -1.68%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "10.00"
$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("E28").Value = "  +2.46%  "

$ws.Range("E29").Value = "  +2.08%  "

$ws.Range("D30").Value = "603.69"
$ws.Range("E30").Value = "  +8.34%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("E32").Value = "  +3.97%  "

$ws.Range("E33").Value = "  +3.83%  "

$ws.Range("E34").Value = "  +5.57%  "

$ws.Range("D35").Value = "0.132"
$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("E36").Value = "  +4.86%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "162.64"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").Value = "20.07"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").Value = "0.382"
$ws.Range("E40").Value = "  +3.03%  "

$ws.Range("D41").Value = "1.92"
$ws.Range("E41").Value = "  +2.43%  "

$ws.Range("D42").Value = "5.46"
$ws.Range("E42").Value = "  +2.48%  "

$ws.Range("D43").Value = "2.69"
$ws.Range("E43").Value = "  +3.30%  "

$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0317"
$ws.Range("E45").Value = "  -4.29%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").Value = "158.57"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "3.95"
$ws.Range("E48").Value = "  +5.63%  "

$ws.Range("E49").Value = "  +6.47%  "

$ws.Range("D50").Value = "0.610"
$ws.Range("E50").Value = "  +7.84%  "

$ws.Range("D51").Value = "22.17"
$ws.Range("E51").Value = "  +0.29%  "
